# Remove footer from slides
#
# The authored change removes the "Footer Placeholder" shape (the
# https://github.com/StoneyJackson/git-ccscne-2016 footer text) from every
# slide in the deck, and clears the inherited footer text on the slide
# master / slide layouts (the placeholder shape itself stays there, only
# its text run is cleared).

$p = $ppt.ActivePresentation

# --- 1. Turn off the footer on every slide -------------------------------
# Setting HeadersFooters.Footer.Visible = $false removes the per-slide
# "Footer Placeholder" shape entirely, which is the core of this edit.
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    try {
        $slide.HeadersFooters.Footer.Visible = $false
    } catch {
    }
}

# --- 2. Clear the footer text inherited from the master / layouts --------
# ppPlaceholderFooter = 15, msoPlaceholder = 14
function Clear-FooterPlaceholderText($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $shp = $shapes.Item($j)
        if ($shp.Type -eq 14 -and $shp.PlaceholderFormat.Type -eq 15) {
            $shp.TextFrame.TextRange.Text = ""
        }
    }
}

$master = $p.SlideMaster
Clear-FooterPlaceholderText $master.Shapes

for ($k = 1; $k -le $master.CustomLayouts.Count; $k++) {
    $layout = $master.CustomLayouts.Item($k)
    Clear-FooterPlaceholderText $layout.Shapes
}
